# Adds the three new Polish (PEF.PL) document type / process identifier rows
# (rows 67-69) to the "Document Type" sheet, matching the PEPPOL Code Lists
# "Document types v6 draft" workbook update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 67: PEF.PL Accounting Note v1
# ---------------------------------------------------------------------
$ws.Range("A66:E66").Copy()
$ws.Range("A67:E67").PasteSpecial(-4104)  # xlPasteAll
$excel.CutCopyMode = $false

$ws.Range("A67").Value = "PEF.PL Accounting Note v1"
$ws.Range("B67").Value = "busdox-docid-qns"
$ws.Range("B67").HorizontalAlignment = -4131  # xlLeft
$ws.Range("C67").Value = "urn:oasis:names:specification:ubl:schema:xsd:CreditNote-2::CreditNote##urn:fdc:www.efaktura.gov.pl:ver1.0:trns:account_corr:ver1.0::2.1"
$ws.Range("D67").Value = 6
$ws.Range("E67").Formula = "=FALSE"

# ---------------------------------------------------------------------
# Row 68: PEF.PL Correcting Invoice v1
# ---------------------------------------------------------------------
$ws.Range("A60:E60").Copy()
$ws.Range("A68:E68").PasteSpecial(-4104)  # xlPasteAll
$excel.CutCopyMode = $false

$ws.Range("A68").Value = "PEF.PL Correcting Invoice v1"
$ws.Range("B68").Value = "busdox-docid-qns"
$ws.Range("B68").HorizontalAlignment = -4131  # xlLeft
$ws.Range("C68").Value = "urn:oasis:names:specification:ubl:schema:xsd:CreditNote-2::CreditNote##urn:cen.eu:en16931:2017#compliant#urn:fdc:peppol.eu:2017:poacc:billing:3.0#extended#urn:fdc:www.efaktura.gov.pl:ver1.0::2.1"
$ws.Range("D68").Value = 6
$ws.Range("E68").Formula = "=FALSE"

# ---------------------------------------------------------------------
# Row 69: PEF.PL Receipt Advice v1
# ---------------------------------------------------------------------
$ws.Range("A60:E60").Copy()
$ws.Range("A69:E69").PasteSpecial(-4104)  # xlPasteAll
$excel.CutCopyMode = $false

$ws.Range("A69").Value = "PEF.PL Receipt Advice v1"
$ws.Range("B69").Value = "busdox-docid-qns"
$ws.Range("B69").HorizontalAlignment = -4131  # xlLeft
$ws.Range("C69").Value = "urn:oasis:names:specification:ubl:schema:xsd:ReceiptAdvice-2::ReceiptAdvice##urn:fdc:www.efaktura.gov.pl:ver1.0:trns:receipt_advice:ver1.0::2.1"
$ws.Range("D69").Value = 6
$ws.Range("E69").Formula = "=FALSE"

# ---------------------------------------------------------------------
# Final selection / scroll state (matches the post-edit cursor position)
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 56
$ws.Range("A71").Select()

Write-Host "Added rows 67-69 (PEF.PL document types)"
